$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new data
$ws.Range("A2").Value = "ROMERO GUSTAVO"
$ws.Range("B2").Value = "'+5491544735404"
$ws.Range("C2").Value = 41000
$ws.Range("D2").Value = "Sent"

# Delete rows 3 and 4 (no longer needed)
$ws.Rows("3:4").Delete()
